$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "24.489.23"
Set-TextValue "E2" "  -1.81%  "
Set-TextValue "D3" "1.670.20"
Set-TextValue "E3" "  -1.99%  "
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "312.64"
Set-TextValue "E5" "  -1.05%  "
Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  -0.05%  "
Set-TextValue "D7" "0.3906"
Set-TextValue "E7" "  -3.87%  "
Set-TextValue "D8" "0.3933"
Set-TextValue "E8" "  -3.32%  "
Set-TextValue "D9" "1.003"
Set-TextValue "E9" "  +0.08%  "
Set-TextValue "D10" "51.41"
Set-TextValue "E10" "  -4.29%  "
Set-TextValue "D11" "1.398"
Set-TextValue "E11" "  -4.95%  "
Set-TextValue "D12" "0.08632"
Set-TextValue "E12" "  -2.11%  "
Set-TextValue "D13" "25.23"
Set-TextValue "E13" "  -2.45%  "
Set-TextValue "D14" "7.304"
Set-TextValue "E14" "  -2.81%  "
Set-TextValue "D15" "0.00001319"
Set-TextValue "E15" "  -2.59%  "
Set-TextValue "D16" "7.701"
Set-TextValue "E16" "  -4.38%  "
Set-TextValue "D17" "1.674.41"
Set-TextValue "E17" "  -3.96%  "
Set-TextValue "D18" "93.31"
Set-TextValue "E18" "  -3.39%  "
Set-TextValue "D19" "0.07013"
Set-TextValue "E19" "  -2.19%  "
Set-TextValue "D20" "20.87"
Set-TextValue "E20" "  -0.68%  "
Set-TextValue "D21" "7.054"
Set-TextValue "E21" "  -2.63%  "
Set-TextValue "E22" "  +0.14%  "
Set-TextValue "E23" "  -4.61%  "
Set-TextValue "D24" "24.493.13"
Set-TextValue "E24" "  -1.80%  "
Set-TextValue "D25" "2.364"
Set-TextValue "E25" "  +1.69%  "
Set-TextValue "D26" "2.736"
Set-TextValue "E26" "  -5.36%  "
Set-TextValue "D28" "5.864"
Set-TextValue "E28" "  -14.27%  "
Set-TextValue "D29" "160.40"
Set-TextValue "E29" "  -2.67%  "
Set-TextValue "D30" "146.95"
Set-TextValue "E30" "  +1.12%  "
Set-TextValue "D31" "8.383"
Set-TextValue "E31" "  +1.52%  "
Set-TextValue "D32" "2.506"
Set-TextValue "E32" "  +10.32%  "
Set-TextValue "D33" "1.852.92"
Set-TextValue "E33" "  -1.94%  "
Set-TextValue "D34" "0.08339"
Set-TextValue "E34" "  -5.26%  "
Set-TextValue "D35" "7.006"
Set-TextValue "E35" "  -4.36%  "
Set-TextValue "D36" "0.03024"
Set-TextValue "E36" "  -5.67%  "
Set-TextValue "D37" "0.2808"
Set-TextValue "E37" "  -1.52%  "
Set-TextValue "D38" "0.9881"
Set-TextValue "E38" "  -2.85%  "
Set-TextValue "D39" "0.09462"
Set-TextValue "E39" "  +0.35%  "
Set-TextValue "D40" "1.513"
Set-TextValue "E40" "  +2.87%  "
Set-TextValue "E41" "  -5.67%  "
Set-TextValue "D42" "0.7884"
Set-TextValue "E42" "  -6.91%  "
Set-TextValue "D43" "13.57"
Set-TextValue "E43" "  -3.58%  "
Set-TextValue "D44" "16.46"
Set-TextValue "E44" "  -7.31%  "
Set-TextValue "D45" "0.7114"
Set-TextValue "E45" "  -4.42%  "
Set-TextValue "D46" "2.545"
Set-TextValue "E46" "  -6.67%  "
Set-TextValue "D47" "4.171"
Set-TextValue "E47" "  -1.53%  "
Set-TextValue "D48" "0.08618"
Set-TextValue "E48" "  +3.10%  "
Set-TextValue "D49" "1.001"
Set-TextValue "E49" "  -0.02%  "
Set-TextValue "D50" "1.321"
Set-TextValue "E50" "  -5.20%  "
Set-TextValue "D51" "137.35"
Set-TextValue "E51" "  -3.47%  "
